# Autogenerated on Thu Mar 26 2015 18:06:15 GMT+0000 (Coordinated Universal Time)
#
# Adds an "Others:" block of additional MSME support organisations (with
# their URLs) underneath the existing "Source:" block on the Summary sheet,
# and duplicates the trailing "SCCSME" citation label/source lines at the
# bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old "Source:" block (rows 44-50) had a hyperlink on the URL line and a
# long-form citation underneath "SCCSME". Both the hyperlink and the long
# citation text go away in the new layout, so drop the hyperlink first and
# then clear out the old rows before rebuilding the block from scratch.
$ws.Hyperlinks.Delete()
$ws.Rows("44:50").Delete()

# Make room for the new content: the rebuilt "Source:" / "Others:" block now
# runs from row 44 through row 65 (one blank spacer row between every pair
# of lines, matching the style already used at row 43 "Source:" / old row 46).
$ws.Rows("44:65").Insert()

# -- Source: (restored, now without the hyperlink) --
$ws.Range("A45").Value = "Support and Consultation Centre for SMEs. Trade and Industry Department of the Hong Kong Government."
$ws.Range("A47").Value = "http://www.success.tid.gov.hk/english/lin_sup_org/gov_dep/service_detail_6863.html"

# -- Others: --
$ws.Range("A49").Value = "Others:"
$ws.Range("A51").Value = "Hong Kong Small and Medium Enterprises General Association"
$ws.Range("A53").Value = "http://www.hksmega.org/english/index.htm"
$ws.Range("A55").Value = "Support for Local Enterprises & SMEs, Government of Hong Kong"
$ws.Range("A57").Value = "http://www.gov.hk/en/business/supportenterprises/localenterprises/"
$ws.Range("A59").Value = "Hong Kong Trade Development Council (HKTDC)"
$ws.Range("A61").Value = "http://www.hktdc.com/mis/ahktdc/en/s/abt-hktdc-about.html"
$ws.Range("A63").Value = "Trade and Industry Department - Support to Small and Medium Enterprises"
$ws.Range("A65").Value = "https://www.tid.gov.hk/english/smes_industry/smes/smes_content.html"

# Rows 66-67 are left untouched (no content, no formatting) so they don't
# appear in the saved sheet at all; the trailing "SCCSME" label/source pair
# lands two rows further down, at 68-69.
$ws.Rows("66:67").ClearContents()
$ws.Rows("66:67").ClearFormats()

$ws.Range("A68").Value = "SCCSME"
$ws.Range("A68").Font.Bold = $true

$ws.Range("A69").Value = "SCCSME"
